# Simula la corrida de "Actualizar" (02-05-2021 16-23-09): refresca la
# marca de tiempo del ultimo batch (filas 660:673) y agrega un nuevo batch
# de 14 filas (674:687) con el mismo patron Nombre/URL/Disponibilidad/Fecha.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Las 14 filas del batch anterior quedan con una marca de tiempo
#        ligeramente distinta (mismo instante de "hoy", recalculado) ---
$ws.Range("D660:D673").Value = 44232.66163265047

# --- 2) Nuevo batch: filas 674-687 ---
$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")
$urls  = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)

$startRow = 674
$newDate = 44232.68268269519

for ($i = 0; $i -lt 14; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $names[$i]
    $ws.Cells.Item($r, 2).Value = $urls[$i]
    $ws.Cells.Item($r, 3).Value = "Disponible"
    $ws.Cells.Item($r, 4).Value = $newDate
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

# Hipervinculos para la columna B del nuevo batch (se crean despues de
# escribir todos los valores para que los Id de relacion salgan
# consecutivos, igual que en el archivo original).
for ($i = 0; $i -lt 14; $i++) {
    $r = $startRow + $i
    $cell = $ws.Cells.Item($r, 2)
    if ($urls[$i] -eq "https://ide.dataintelligence-group.com/mapstore/#/") {
        $ws.Hyperlinks.Add($cell, "https://ide.dataintelligence-group.com/mapstore/", "/")
    } else {
        $ws.Hyperlinks.Add($cell, $urls[$i])
    }
}

# Restaura el estilo "Hyperlink" (subrayado/tema) en la columna B: Add()
# aplica su propio formato interno; lo igualamos al del resto de la hoja.
for ($i = 0; $i -lt 14; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Style = "Hyperlink"
}
